$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB -> BNB
$ws.Range("D2").Value = "'328.46"
$ws.Range("E2").Value = "'1.33%"
$ws.Range("G2").Value = "'7"

# Row 3: OKB -> OKB
$ws.Range("D3").Value = "'43.93"
$ws.Range("E3").Value = "'-1.09%"
$ws.Range("G3").Value = "'7"

# Row 4: HuobiToken -> HuobiToken
$ws.Range("D4").Value = "'5.516"
$ws.Range("E4").Value = "'0.56%"
$ws.Range("G4").Value = "'7"

# Row 5: Cronos -> Cronos
$ws.Range("E5").Value = "'-0.13%"
$ws.Range("G5").Value = "'7"

# Row 6: FTXToken -> FTXToken
$ws.Range("D6").Value = "'1.994"
$ws.Range("E6").Value = "'4.60%"
$ws.Range("G6").Value = "'7"

# Row 7: GateToken -> GateToken
$ws.Range("D7").Value = "'4.358"
$ws.Range("E7").Value = "'1.71%"
$ws.Range("G7").Value = "'7"

# Row 8: BTSEToken -> BTSEToken
$ws.Range("E8").Value = "'-5.48%"
$ws.Range("G8").Value = "'7"

# Row 9: MXToken -> MXToken
$ws.Range("D9").Value = "'0.9493"
$ws.Range("E9").Value = "'1.09%"
$ws.Range("G9").Value = "'7"

# Row 10: LiechtensteinCryptoassetsExchange -> LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.1132"
$ws.Range("E10").Value = "'-4.30%"
$ws.Range("G10").Value = "'7"

# Row 11: WazirX -> WazirX
$ws.Range("E11").Value = "'0.54%"
$ws.Range("G11").Value = "'7"

# Row 12: MCDex -> MCDex
$ws.Range("D12").Value = "'10.72"
$ws.Range("E12").Value = "'25.98%"
$ws.Range("G12").Value = "'7"

# Row 13: MandalaExchangeToken -> MandalaExchangeToken
$ws.Range("D13").Value = "'0.09914"
$ws.Range("E13").Value = "'-0.47%"
$ws.Range("G13").Value = "'7"

# Row 14: BitrueCoin -> BitrueCoin
$ws.Range("D14").Value = "'0.04674"
$ws.Range("E14").Value = "'7.02%"
$ws.Range("G14").Value = "'7"

# Row 15: BitMartToken -> BitMartToken
$ws.Range("D15").Value = "'0.1066"
$ws.Range("E15").Value = "'0.20%"
$ws.Range("G15").Value = "'7"

# Row 16: BitForexToken -> BitForexToken
$ws.Range("D16").Value = "'0.001276"
$ws.Range("E16").Value = "'-0.78%"
$ws.Range("G16").Value = "'7"

# Row 17: CoinExToken -> CoinExToken
$ws.Range("D17").Value = "'0.04075"
$ws.Range("E17").Value = "'-4.14%"
$ws.Range("G17").Value = "'7"

# Row 18: TigerCash -> TigerCash
$ws.Range("D18").Value = "'0.005973"
$ws.Range("E18").Value = "'1.68%"
$ws.Range("G18").Value = "'7"

# Row 19: HotbitToken -> LEO
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.359"
$ws.Range("E19").Value = "'-6.50%"
$ws.Range("G19").Value = "'7"

# Row 20: LEO -> BitpandaEcosystemToken
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3475"
$ws.Range("E20").Value = "'-0.23%"
$ws.Range("G20").Value = "'7"

# Row 21: BitpandaEcosystemToken -> ProBitToken
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1420"
$ws.Range("E21").Value = "'3.54%"
$ws.Range("G21").Value = "'7"

# Row 22: ProBitToken -> ZBToken
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2545"
$ws.Range("E22").Value = "'0.49%"
$ws.Range("G22").Value = "'7"

# Row 23: ZBToken -> BitKan
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001263"
$ws.Range("E23").Value = "'2.13%"
$ws.Range("G23").Value = "'7"

# Row 24: BitKan -> HotbitToken
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004326"
$ws.Range("E24").Value = "'-4.62%"
$ws.Range("G24").Value = "'7"

# Row 25: NitroEx -> NitroEx
$ws.Range("D25").Value = "'0.0001201"
$ws.Range("E25").Value = "'-0.14%"
$ws.Range("G25").Value = "'7"

# Row 26: UpBots -> UpBots
$ws.Range("D26").Value = "'0.0003746"
$ws.Range("E26").Value = "'-6.27%"
$ws.Range("G26").Value = "'7"

# Row 27: Spectre.aiUtilityToken -> Spectre.aiUtilityToken
$ws.Range("G27").Value = "'7"

# Row 28: LegolasExchange -> LegolasExchange
$ws.Range("G28").Value = "'7"

# Row 29: BitZToken -> BitZToken
$ws.Range("G29").Value = "'7"

# Row 30: Birake -> Birake
$ws.Range("G30").Value = "'7"

# Row 31: NashExchange -> NashExchange
$ws.Range("G31").Value = "'7"

# Row 32: AAXToken -> AAXToken
$ws.Range("G32").Value = "'7"

# Row 33: CenX -> CenX
$ws.Range("G33").Value = "'7"

# Row 34: BNIXToken -> BNIXToken
$ws.Range("G34").Value = "'7"

# Row 35: Polkally -> Polkally
$ws.Range("G35").Value = "'7"

# Row 36: Charli3 -> Charli3
$ws.Range("G36").Value = "'7"

# Row 37: BlubitexToken -> BlubitexToken
$ws.Range("G37").Value = "'7"

# Row 38: One -> One
$ws.Range("D38").Value = "'0.02587"
$ws.Range("E38").Value = "'-1.80%"
$ws.Range("G38").Value = "'7"

# Row 39: IDEX -> IDEX
$ws.Range("D39").Value = "'0.05641"
$ws.Range("E39").Value = "'3.05%"
$ws.Range("G39").Value = "'7"

# Row 40: KickToken -> KickToken
$ws.Range("D40").Value = "'0.007560"
$ws.Range("E40").Value = "'-1.42%"
$ws.Range("G40").Value = "'7"

# Row 41: BKEXToken -> BKEXToken
$ws.Range("D41").Value = "'0.1396"
$ws.Range("E41").Value = "'0.39%"
$ws.Range("G41").Value = "'7"

# Row 42: Dexo -> Dexo
$ws.Range("D42").Value = "'0.007503"
$ws.Range("E42").Value = "'5.09%"
$ws.Range("G42").Value = "'7"

# Row 43: CEJI -> CEJI
$ws.Range("D43").Value = "'0.002016"
$ws.Range("E43").Value = "'-3.16%"
$ws.Range("G43").Value = "'7"

# Row 44: LocalTraders -> LocalTraders
$ws.Range("E44").Value = "'-1.11%"
$ws.Range("G44").Value = "'7"

# Row 45: CoinLion -> CoinLion
$ws.Range("D45").Value = "'0.00007102"
$ws.Range("E45").Value = "'-0.28%"
$ws.Range("G45").Value = "'7"

# Row 46: Kangarootoken -> Kangarootoken
$ws.Range("E46").Value = "'-0.15%"
$ws.Range("G46").Value = "'7"

# Row 47: CoinbaseStockToken -> CoinbaseStockToken
$ws.Range("E47").Value = "'55.16%"
$ws.Range("G47").Value = "'7"

# Row 48: BOLO -> BOLO
$ws.Range("D48").Value = "'0.003706"
$ws.Range("E48").Value = "'-0.34%"
$ws.Range("G48").Value = "'7"

# Row 49: CryptobidCoin -> CryptobidCoin
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.15%"
$ws.Range("G49").Value = "'7"

# Row 50: SpecialPowerGold -> SpecialPowerGold
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("G50").Value = "'7"

# Row 51: DigiFinexToken -> DigiFinexToken
$ws.Range("G51").Value = "'7"
